# Refresh the cryptos price/volume snapshot (commit: "Updated cryptos list
# on Sat Nov 18 05:24:50 UTC 2023 with GitHub Actions").
#
# Every data cell on the sheet is stored as text (even price columns that
# look numeric, e.g. "241.03"), so numeric-looking replacement values are
# written with a leading apostrophe to keep Excel from auto-converting them
# to real numbers (which would silently drop trailing zeros, e.g. "1.00" -> 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.383.57"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "1.932.52"
$ws.Range("E3").Value = "  -2.49%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'241.03"
$ws.Range("E5").Value = "  -1.33%  "
$ws.Range("D6").Value = "'0.607"
$ws.Range("E6").Value = "  -3.38%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'56.84"
$ws.Range("E8").Value = "  -4.54%  "
$ws.Range("D9").Value = "'0.358"
$ws.Range("E9").Value = "  -4.58%  "
$ws.Range("D10").Value = "'0.0832"
$ws.Range("E10").Value = "  +1.47%  "
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("D12").Value = "2.217.03"
$ws.Range("E12").Value = "  -2.35%  "
$ws.Range("D13").Value = "'0.801"
$ws.Range("E13").Value = "  -7.27%  "
$ws.Range("D14").Value = "'20.95"
$ws.Range("E14").Value = "  -12.04%  "
$ws.Range("D15").Value = "'13.36"
$ws.Range("E15").Value = "  -4.59%  "
$ws.Range("D16").Value = "'5.14"
$ws.Range("E16").Value = "  -6.09%  "
$ws.Range("D17").Value = "1.933.02"
$ws.Range("E17").Value = "  -2.40%  "
$ws.Range("D18").Value = "36.274.57"
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("B19").Value = "Litecoin"
$ws.Range("C19").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D19").Value = "'68.95"
$ws.Range("E19").Value = "  -1.23%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0860"
$ws.Range("E20").Value = "  -0.59%  "
$ws.Range("D21").Value = "'227.39"
$ws.Range("E21").Value = "  -2.92%  "
$ws.Range("D22").Value = "'4.95"
$ws.Range("E22").Value = "  -6.76%  "
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("E24").Value = "  -10.14%  "
$ws.Range("E25").Value = "  -1.69%  "
$ws.Range("D26").Value = "'9.26"
$ws.Range("E26").Value = "  -7.77%  "
$ws.Range("D27").Value = "'160.50"
$ws.Range("E27").Value = "  -1.04%  "
$ws.Range("E28").Value = "  -2.60%  "
$ws.Range("D29").Value = "'19.15"
$ws.Range("E29").Value = "  -3.37%  "
$ws.Range("E30").Value = "  -2.29%  "
$ws.Range("D31").Value = "'1.11"
$ws.Range("E31").Value = "  -6.67%  "
$ws.Range("D32").Value = "'4.55"
$ws.Range("E32").Value = "  -7.27%  "
$ws.Range("D33").Value = "'0.0626"
$ws.Range("E33").Value = "  -0.34%  "
$ws.Range("D34").Value = "'4.14"
$ws.Range("E34").Value = "  -6.77%  "
$ws.Range("B35").Value = "BinanceUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("B36").Value = "THORChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D36").Value = "'6.10"
$ws.Range("E36").Value = "  -2.09%  "
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("D38").Value = "'2.13"
$ws.Range("E38").Value = "  -6.46%  "
$ws.Range("E39").Value = "  -1.76%  "
$ws.Range("D40").Value = "'0.0971"
$ws.Range("E40").Value = "  +0.31%  "
$ws.Range("E41").Value = "  -1.47%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "'0.0208"
$ws.Range("E42").Value = "  -2.76%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'1.15"
$ws.Range("E43").Value = "  -7.74%  "
$ws.Range("D44").Value = "'15.42"
$ws.Range("E44").Value = "  -5.24%  "
$ws.Range("D45").Value = "1.334.04"
$ws.Range("E45").Value = "  -2.80%  "
$ws.Range("D46").Value = "'1.01"
$ws.Range("E46").Value = "  -7.41%  "
$ws.Range("D47").Value = "'86.55"
$ws.Range("E47").Value = "  -6.58%  "
$ws.Range("D48").Value = "'7.07"
$ws.Range("E48").Value = "  -5.73%  "
$ws.Range("D49").Value = "'2.82"
$ws.Range("E49").Value = "  -0.88%  "
$ws.Range("D50").Value = "'43.86"
$ws.Range("E50").Value = "  -3.30%  "
$ws.Range("D51").Value = "2.108.99"
$ws.Range("E51").Value = "  -2.49%  "
